# Generate Report for Handoff
# Update "Latest Handoff Datetime" (column G) for the
# 85049ea7-dd01-4e3e-9940-3caa4db5c7a8 entry (row 4) on both the
# zh-cn and de-de detail sheets to reflect a freshly generated handoff.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("G4").Value = "2016-07-26 07:27:29"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("G4").Value = "2016-07-26 07:27:38"
